# Applies the "Add files via upload" update to InformacjeOPrzeniesieniach.xlsx
#  1. Fix a room-number typo in an existing "Oddzialy" row (sala: 38 -> sala: 37)
#  2. Insert two new transfer rows (for 19.12.2025 lessons 6 and 7) before the
#     existing last data row, which shifts it from row 16 down to row 18.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Oddziały")

# 1) Correct the room number for the Najwer/Informatyka transfer row
$ws.Cells.Item(12, 2).Value = "18.12.2025, 6, 12:25-13:10, sala: 37"

# 2) Insert two blank rows above the current last row (row 16), pushing the
#    "Jarek Zbigniew" row down to row 18.
$ws.Range("A16:A17").EntireRow.Insert()

# New row 16: Misiąg Anna / 3TH|JA1 / Język niemiecki
$ws.Cells.Item(16, 1).Value = "19.12.2025, 6, 12:25-13:10, sala: 44"
$ws.Cells.Item(16, 2).Value = "19.12.2025, 6, 12:25-13:10, sala: 42"
$ws.Cells.Item(16, 3).Value = "Misiąg Anna"
$ws.Cells.Item(16, 4).Value = "-"
$ws.Cells.Item(16, 5).Value = "3TH|JA1"
$ws.Cells.Item(16, 6).Value = "Język niemiecki"
$ws.Cells.Item(16, 7).Value = ""

# New row 17: Biczysko Wojciech / 3TH / Fizyka
$ws.Cells.Item(17, 1).Value = "19.12.2025, 7, 13:15-14:00, sala: 22"
$ws.Cells.Item(17, 2).Value = "19.12.2025, 7, 13:15-14:00, sala: 42"
$ws.Cells.Item(17, 3).Value = "Biczysko Wojciech"
$ws.Cells.Item(17, 4).Value = "-"
$ws.Cells.Item(17, 5).Value = "3TH"
$ws.Cells.Item(17, 6).Value = "Fizyka"
$ws.Cells.Item(17, 7).Value = ""
